$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Header row tweak: wrap text on the "Solution" / "Phase" header cells
# ---------------------------------------------------------------------------
$ws.Range("C2:D2").WrapText = $true

# ---------------------------------------------------------------------------
# 2) Apply the same look as the existing data rows (3-7) to row 8, which was
#    a blank styled row and now becomes a normal data row.
# ---------------------------------------------------------------------------
$ws.Range("A7:E7").Copy()
$ws.Range("A8:E8").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 3) New reflection entries (rows 8-12), "Management" phase
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "Member don't care about Risks in project and don't update Risk Category"
$ws.Range("C8").Value = "Follow risk plan, plan about interation for risk"
$ws.Range("D8").Value = "Knowing more about manage risk better"
$ws.Range("E8").Value = "Management"

$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "Too much and more difficult to measurement"
$ws.Range("C9").Value = "Research about measurement, implement Goal-Question-Metric"
$ws.Range("D9").Value = "Knowing more defenite about metrics and how to get it"
$ws.Range("E9").Value = "Management"

$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "Project difficult to control and monitoring"
$ws.Range("C10").Value = "Plan for detail plan, WBS, implement tracking and monitoring through measurement about schedule devition metric"
$ws.Range("E10").Value = "Management"

$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "360 review is not good conduct"
$ws.Range("C11").Value = "Require team member write reflection base on 360 reivew"
$ws.Range("D11").Value = "Knowing about management and communicate between team member"
$ws.Range("E11").Value = "Management"

$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "Team member is not complete work on time"
$ws.Range("C12").Value = "Re-estimate, and evaluate effort of team member"
$ws.Range("D12").Value = "Conduct measurement about productivity"
$ws.Range("E12").Value = "Management"

# ---------------------------------------------------------------------------
# 4) Row heights for the newly-filled rows
# ---------------------------------------------------------------------------
$ws.Rows("8").RowHeight = 30
$ws.Rows("9").RowHeight = 30
$ws.Rows("10").RowHeight = 45
$ws.Rows("11").RowHeight = 30
$ws.Rows("12").RowHeight = 30

# ---------------------------------------------------------------------------
# 5) Column A (the "No.") keeps the centred / wrapped look all the way to
#    row 18, and B:E keep a numbered sequence / border but without the
#    wrap+center treatment used in the main table body.
# ---------------------------------------------------------------------------
$ws.Range("A8").Copy()
$ws.Range("A9:A18").PasteSpecial(-4122)

$ws.Range("A13").Value = 11
$ws.Range("A14").Value = 12
$ws.Range("A15").Value = 13
$ws.Range("A16").Value = 14
$ws.Range("A17").Value = 15
$ws.Range("A18").Value = 16

# B, C, D, E for rows 9-19 need a bordered Times-New-Roman style matching the
# rest of the table. Start from the wrapped data-row look (so the correct
# font gets reused) then switch off wrapping where required.
$ws.Range("B3").Copy()
$ws.Range("B9:E19").PasteSpecial(-4122)
$ws.Range("C9:D19").WrapText = $true
$ws.Range("B9:B19").WrapText = $false
$ws.Range("E9:E19").WrapText = $false
$ws.Range("A19").WrapText = $false

Write-Host "done"
